# Applies the cell-value fills described in the commit "filled in og data vals".
# For each metric block (Precision/Recall/Reciprocal rank/DCG/ERR under columns
# C-G, L-P, U-Y, AD-AH, AM-AQ, AV-AZ) the previously blank or placeholder score rows
# are populated with their real computed values, and a couple of already-populated
# rows are corrected to the right numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 0.75
$ws.Range("D4").Value = 0.68181818181800002
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 22.399506517700001
$ws.Range("G4").Value = 1.4634587241399999
$ws.Range("L4").Value = 0.55000000000000004
$ws.Range("M4").Value = 0.5
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 15.0501008584
$ws.Range("P4").Value = 1.4634528678300001
$ws.Range("AD4").Value = 0.95
$ws.Range("AE4").Value = 0.86363636363600005
$ws.Range("AF4").Value = 1
$ws.Range("AG4").Value = 30.053242021599999
$ws.Range("AH4").Value = 1.46345872806
$ws.Range("AM4").Value = 0.25
$ws.Range("AN4").Value = 0.22727272727299999
$ws.Range("AO4").Value = 0.25
$ws.Range("AP4").Value = 1.53157201951
$ws.Range("AQ4").Value = 0.018221378548699999
$ws.Range("AV4").Value = 0.15
$ws.Range("AW4").Value = 0.13636363636400001
$ws.Range("AX4").Value = 0.125
$ws.Range("AY4").Value = 0.82142290159599995
$ws.Range("AZ4").Value = 0.0080868675595200008

# Row 18
$ws.Range("C18").Value = 0.75
$ws.Range("D18").Value = 0.68181818181800002
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 22.399506517700001
$ws.Range("G18").Value = 1.4634587241399999
$ws.Range("L18").Value = 0.7
$ws.Range("M18").Value = 0.63636363636399995
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 8.2244342271100006
$ws.Range("P18").Value = 0.38266054708199998
$ws.Range("AD18").Value = 1
$ws.Range("AE18").Value = 0.90909090909099999
$ws.Range("AF18").Value = 1
$ws.Range("AG18").Value = 7.0402683819199998
$ws.Range("AH18").Value = 0.100718798328
$ws.Range("AM18").Value = 0.57894736842100003
$ws.Range("AN18").Value = 0.5
$ws.Range("AO18").Value = 0.33333333333300003
$ws.Range("AP18").Value = 13.478512328500001
$ws.Range("AQ18").Value = 0.34180626494100003
$ws.Range("AV18").Value = 0.2
$ws.Range("AW18").Value = 0.181818181818
$ws.Range("AX18").Value = 0.111111111111
$ws.Range("AY18").Value = 1.0424588298199999
$ws.Range("AZ18").Value = 0.0092154605551899998

# Row 33
$ws.Range("C33").Value = 0.75
$ws.Range("D33").Value = 0.68181818181800002
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 22.399506517700001
$ws.Range("G33").Value = 1.4634587241399999
$ws.Range("L33").Value = 0.25
$ws.Range("M33").Value = 0.22727272727299999
$ws.Range("N33").Value = 1
$ws.Range("O33").Value = 9.4721781331700008
$ws.Range("P33").Value = 0.45346406300999997
$ws.Range("U33").Value = 0.3
$ws.Range("V33").Value = 0.31276583899999999
$ws.Range("W33").Value = 0.5
$ws.Range("X33").Value = 4.5673298200000003
$ws.Range("Y33").Value = 0.56203285319999996
$ws.Range("AD33").Value = 0.7
$ws.Range("AE33").Value = 0.63636363636399995
$ws.Range("AF33").Value = 1
$ws.Range("AG33").Value = 21.391061511299998
$ws.Range("AH33").Value = 1.46345872806
$ws.Range("AM33").Value = 0.5
$ws.Range("AN33").Value = 0.45454545454500001
$ws.Range("AO33").Value = 0.33333333333300003
$ws.Range("AP33").Value = 5.0566816920299997
$ws.Range("AQ33").Value = 0.32370127735600002
$ws.Range("AV33").Value = 0.05
$ws.Range("AW33").Value = 0.045454545454499999
$ws.Range("AX33").Value = 0.083333333333299994
$ws.Range("AY33").Value = 0.270238154427
$ws.Range("AZ33").Value = 0.0026041666666699998

# Row 41
$ws.Range("C41").Value = 0.75
$ws.Range("D41").Value = 0.68181818181800002
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 22.399506517700001
$ws.Range("G41").Value = 1.4634587241399999

# Row 47
$ws.Range("C47").Value = 0.75
$ws.Range("D47").Value = 0.68181818181800002
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = 22.399506517700001
$ws.Range("G47").Value = 1.4634587241399999
$ws.Range("U47").Value = 0.85
$ws.Range("V47").Value = 0.77272727272700004
$ws.Range("W47").Value = 1
$ws.Range("X47").Value = 25.379035872799999
$ws.Range("Y47").Value = 1.4634587241699999
$ws.Range("AD47").Value = 0.65
$ws.Range("AE47").Value = 0.59090909090900001
$ws.Range("AF47").Value = 1
$ws.Range("AG47").Value = 21.0900315156
$ws.Range("AH47").Value = 1.46345872806
$ws.Range("AM47").Value = 0.95
$ws.Range("AN47").Value = 0.86363636363600005
$ws.Range("AO47").Value = 1
$ws.Range("AP47").Value = 28.134479678000002
$ws.Range("AQ47").Value = 1.46345872806
$ws.Range("AV47").Value = 0.05
$ws.Range("AW47").Value = 0.045454545454499999
$ws.Range("AX47").Value = 0.083333333333299994
$ws.Range("AY47").Value = 0.270238154427
$ws.Range("AZ47").Value = 0.0026041666666699998

# Row 51
$ws.Range("C51").Value = 0.75
$ws.Range("D51").Value = 0.68181818181800002
$ws.Range("E51").Value = 1
$ws.Range("F51").Value = 22.399506517700001
$ws.Range("G51").Value = 1.4634587241399999

# Row 62
$ws.Range("C62").Value = 0.75
$ws.Range("D62").Value = 0.68181818181800002
$ws.Range("E62").Value = 1
$ws.Range("F62").Value = 22.399506517700001
$ws.Range("G62").Value = 1.4634587241399999
$ws.Range("AD62").Value = 1
$ws.Range("AE62").Value = 0.90909090909099999
$ws.Range("AF62").Value = 1
$ws.Range("AG62").Value = 34.275829057000003
$ws.Range("AH62").Value = 1.46345872806
$ws.Range("AM62").Value = 1
$ws.Range("AN62").Value = 0.90909090909099999
$ws.Range("AO62").Value = 1
$ws.Range("AP62").Value = 32.283595263400002
$ws.Range("AQ62").Value = 1.46345872806
$ws.Range("AV62").Value = 0.8
$ws.Range("AW62").Value = 0.72727272727299996
$ws.Range("AX62").Value = 1
$ws.Range("AY62").Value = 24.2238668516
$ws.Range("AZ62").Value = 1.4634585738500001

# Row 68
$ws.Range("C68").Value = 0.75
$ws.Range("D68").Value = 0.68181818181800002
$ws.Range("E68").Value = 1
$ws.Range("F68").Value = 22.399506517700001
$ws.Range("G68").Value = 1.4634587241399999

# Row 76
$ws.Range("C76").Value = 0.75
$ws.Range("D76").Value = 0.68181818181800002
$ws.Range("E76").Value = 1
$ws.Range("F76").Value = 22.399506517700001
$ws.Range("G76").Value = 1.4634587241399999
$ws.Range("AD76").Value = 1
$ws.Range("AE76").Value = 0.90909090909099999
$ws.Range("AF76").Value = 1
$ws.Range("AG76").Value = 35.201341909600004
$ws.Range("AH76").Value = 1.46345872806
$ws.Range("AM76").Value = 0.25
$ws.Range("AN76").Value = 0.22727272727299999
$ws.Range("AO76").Value = 0.2
$ws.Range("AP76").Value = 1.3211226490300001
$ws.Range("AQ76").Value = 0.0126767406359
$ws.Range("AV76").Value = 1
$ws.Range("AW76").Value = 0.90909090909099999
$ws.Range("AX76").Value = 1
$ws.Range("AY76").Value = 10.5056230497
$ws.Range("AZ76").Value = 0.52811091561800005

# Row 77
$ws.Range("C77").Value = 0.75
$ws.Range("D77").Value = 0.68181818181800002
$ws.Range("E77").Value = 1
$ws.Range("F77").Value = 22.399506517700001
$ws.Range("G77").Value = 1.4634587241399999

# Update the sheet selection to match the saved view state
$ws.Range("AV76:AZ76").Select()
